$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G6").Value = 3
$ws.Range("I6").Value = 2.7
$ws.Range("J6").Value = 4
$ws.Range("L6").Value = 3.6
$ws.Range("O6").Value = 1.62
$ws.Range("P6").Value = 2.2
$ws.Range("Q6").Value = 3.1
$ws.Range("R6").Value = 1.36
$ws.Range("S6").Value = 1.67
$ws.Range("T6").Value = 2.1
$ws.Range("W6").Value = 6.5
$ws.Range("X6").Value = 13
$ws.Range("Z6").Value = 34
$ws.Range("AD6").Value = 5.5
$ws.Range("AG6").Value = 6
$ws.Range("AH6").Value = 11
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 29
$ws.Range("AK6").Value = 29
$ws.Range("AN6").Value = 4.75
$ws.Range("AO6").Value = 19
$ws.Range("AP6").Value = 34
$ws.Range("AQ6").Value = 67
$ws.Range("AT6").Value = 2.1
$ws.Range("AW6").Value = 4.5
$ws.Range("AX6").Value = 17
